$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.982.99"
$ws.Range("E2").Value = "  +1.67%  "

$ws.Range("D3").Value = "2.249.66"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.81"
$ws.Range("E5").Value = "  +2.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.87"
$ws.Range("E6").Value = "  +3.38%  "

$ws.Range("E7").Value = "  +1.21%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  +1.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.17"
$ws.Range("E10").Value = "  +5.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.84"
$ws.Range("E11").Value = "  +6.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0792"
$ws.Range("E12").Value = "  +1.02%  "

$ws.Range("E13").Value = "  +2.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.57"
$ws.Range("E14").Value = "  +0.66%  "

$ws.Range("D15").Value = "2.590.25"
$ws.Range("E15").Value = "  +0.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.15"
$ws.Range("E16").Value = "  +1.93%  "

$ws.Range("D17").Value = "2.239.06"
$ws.Range("E17").Value = "  +0.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.751"
$ws.Range("E18").Value = "  +1.88%  "

$ws.Range("D19").Value = "40.871.24"
$ws.Range("E19").Value = "  +1.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.82"
$ws.Range("E20").Value = "  +3.25%  "

$ws.Range("D21").Value = "0.0₃0902"
$ws.Range("E21").Value = "  +1.15%  "

$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.54"
$ws.Range("E23").Value = "  +1.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.90"
$ws.Range("E24").Value = "  +1.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.56"
$ws.Range("E25").Value = "  +3.43%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.87"
$ws.Range("E27").Value = "  +2.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.99"
$ws.Range("E28").Value = "  +5.02%  "

$ws.Range("E29").Value = "  -1.89%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.51"
$ws.Range("E30").Value = "  +2.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.08"
$ws.Range("E31").Value = "  +1.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.39"
$ws.Range("E32").Value = "  +3.39%  "

$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.10"
$ws.Range("E34").Value = "  +2.82%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.05"
$ws.Range("E35").Value = "  +5.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0728"
$ws.Range("E36").Value = "  +1.23%  "

$ws.Range("E37").Value = "  +7.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.35"
$ws.Range("E38").Value = "  -0.44%  "

$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.46"
$ws.Range("E39").Value = "  +3.47%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.115"
$ws.Range("E40").Value = "  +2.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.79"
$ws.Range("E41").Value = "  +5.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.91"
$ws.Range("E42").Value = "  +1.20%  "

$ws.Range("D43").Value = "2.093.35"
$ws.Range("E43").Value = "  -2.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.79"
$ws.Range("E44").Value = "  +8.16%  "

$ws.Range("E45").Value = "  +2.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.23"
$ws.Range("E46").Value = "  +4.00%  "

$ws.Range("E47").Value = "  +8.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.83"
$ws.Range("E48").Value = "  -14.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.54"
$ws.Range("E49").Value = "  +3.11%  "

$ws.Range("D50").Value = "2.461.85"
$ws.Range("E50").Value = "  +0.57%  "

$ws.Range("E51").Value = "  +3.10%  "
